# 12/02/2019 : Remplacement des gros boutons par des plus petits
# Update row 16 of "BOM_format_RS" sheet: replace SKQGAKE010 button part
# with the smaller PTS810 SJG 250 SMTR LFS button.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BOM_format_RS")

$ws.Range("B16").Value = "Bouton"
$ws.Range("H16").Value = "135-9492"
$ws.Range("C16").Value = "PTS810 SJG 250 SMTR LFS"
$ws.Range("E16").Value = "PTS810 SJG 250 SMTR LFS"
$ws.Range("F16").Value = "PTS810 SJG 250 SMTR LFS"
$ws.Range("I16").Value = 20
$ws.Range("J16").Value = 0.285

# The old hyperlink/reference text in M16 goes away, but the cell keeps the
# Hyperlink style Excel created for it.
$ws.Range("M16").Style = "Hyperlink"
$ws.Range("M16").Value = ""

# Drop the now-obsolete supplemental columns for this row.
$ws.Range("N16").ClearContents()
$ws.Range("P16").ClearContents()
$ws.Range("Q16").ClearContents()
$ws.Range("R16").ClearContents()
$ws.Range("S16").ClearContents()
$ws.Range("T16").ClearContents()

# Selection left where the author's cursor ended up after editing.
$ws.Range("I17").Select()
